# Update benchmark: 2025-11-21 06:39:55 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("D4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("D5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("D6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

$ws.Range("D8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"

$ws.Range("D9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"

$ws.Range("D10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

$ws.Range("D11").Value = "3,04 TL - 6,09 TL - 76,17 TL"

$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"

$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"

$ws.Range("D14").Value = "2.300 TL - 9.500 TL"
$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
